# Apply HLM test-data update: replace "北京" (Beijing) station data with
# "敦煌 " (Dunhuang) data on both sheets, and update the associated readings.

$wb = $excel.ActiveWorkbook

# --- Sheet "erosion" (sheet1) ---
$ws1 = $wb.Worksheets.Item("erosion")

$ws1.Range("A3").Value = 1
$ws1.Range("C2").Value = "敦煌 "
$ws1.Range("C3").Value = "敦煌 "

$ws1.Columns.Item(2).ColumnWidth = 29.7

$ws1.Range("A3").Select()

# --- Sheet "soil" (sheet2) ---
$ws2 = $wb.Worksheets.Item("soil")

$ws2.Range("B2").Value = 11.3
$ws2.Range("C2").Value = 38
$ws2.Range("D2").Value = 48
$ws2.Range("E2").Value = 7
$ws2.Range("F2").Value = 2632
$ws2.Range("G2").Value = 7846
$ws2.Range("H2").Value = "敦煌 "

$ws2.Range("H2").Select()
